$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1189.9286
$ws.Range("I98").Value = 1011.25
$ws.Range("J98").Value = 1636.625
$ws.Range("K98").Value = 1011.25
$ws.Range("L98").Value = 1636.625
$ws.Range("M98").Value = 486.75
$ws.Range("N98").Value = -4632.625

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3572.9473
$ws.Range("I113").Value = 2335.5557
$ws.Range("J113").Value = 4686.6
$ws.Range("K113").Value = 2335.5557
$ws.Range("L113").Value = 4686.6
$ws.Range("M113").Value = 918.4443000000001
$ws.Range("N113").Value = -11194.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1189.9286
$ws.Range("I122").Value = 1011.25
$ws.Range("J122").Value = 1636.625
$ws.Range("K122").Value = 3033.75
$ws.Range("L122").Value = 4909.875
$ws.Range("M122").Value = -583.75
$ws.Range("N122").Value = -9809.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5862.4434
$ws.Range("I32").Value = 5564.2593
$ws.Range("K32").Value = 5564.2593
$ws.Range("M32").Value = -5277.2593

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1677.7667
$ws.Range("I74").Value = 1643.762
$ws.Range("J74").Value = 1757.1111
$ws.Range("K74").Value = 1643.762
$ws.Range("L74").Value = 1757.1111
$ws.Range("M74").Value = -769.7619999999999
$ws.Range("N74").Value = -3505.1111

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1677.7667
$ws.Range("I77").Value = 1643.762
$ws.Range("J77").Value = 1757.1111
$ws.Range("K77").Value = 8218.809999999999
$ws.Range("L77").Value = 8785.5555
$ws.Range("M77").Value = -3850.809999999999
$ws.Range("N77").Value = -17521.5555

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1834.5416
$ws.Range("I102").Value = 1533.1052
$ws.Range("K102").Value = 1533.1052
$ws.Range("M102").Value = 88.89480000000003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H124").Value = 20101.9
$ws.Range("J124").Value = 20101.9
$ws.Range("L124").Value = 20101.9
$ws.Range("N124").Value = -29921.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H48").Value = 234190.8
$ws.Range("J48").Value = 234190.8
$ws.Range("L48").Value = 234190.8
$ws.Range("N48").Value = -235020.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H70").Value = 180000
$ws.Range("J70").Value = 180000
$ws.Range("L70").Value = 180000
$ws.Range("N70").Value = -180586

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H73").Value = 180000
$ws.Range("J73").Value = 180000
$ws.Range("L73").Value = 180000
$ws.Range("N73").Value = -182028

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 9305.951999999999
$ws.Range("I134").Value = 5002.4287
$ws.Range("J134").Value = 11457.714
$ws.Range("K134").Value = 15007.2861
$ws.Range("L134").Value = 34373.142
$ws.Range("M134").Value = -12472.2861
$ws.Range("N134").Value = -39443.142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6291131.5
$ws.Range("I31").Value = 1288.7693
$ws.Range("J31").Value = 23812836
$ws.Range("K31").Value = 1288.7693
$ws.Range("L31").Value = 23812836
$ws.Range("M31").Value = -993.7692999999999
$ws.Range("N31").Value = -23813426

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6291131.5
$ws.Range("I34").Value = 1288.7693
$ws.Range("J34").Value = 23812836
$ws.Range("K34").Value = 1288.7693
$ws.Range("L34").Value = 23812836
$ws.Range("M34").Value = -1086.7693
$ws.Range("N34").Value = -23813240

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3423.5576
$ws.Range("I58").Value = 1541.0278
$ws.Range("J58").Value = 7659.25
$ws.Range("K58").Value = 1541.0278
$ws.Range("L58").Value = 7659.25
$ws.Range("M58").Value = -1338.0278
$ws.Range("N58").Value = -8065.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1537.8422
$ws.Range("I134").Value = 794.2308
$ws.Range("J134").Value = 3149
$ws.Range("K134").Value = 2382.6924
$ws.Range("L134").Value = 9447
$ws.Range("M134").Value = 152.3076000000001
$ws.Range("N134").Value = -14517

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3423.5576
$ws.Range("I136").Value = 1541.0278
$ws.Range("J136").Value = 7659.25
$ws.Range("K136").Value = 4623.0834
$ws.Range("L136").Value = 22977.75
$ws.Range("M136").Value = -2073.0834
$ws.Range("N136").Value = -28077.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1030.7646
$ws.Range("J97").Value = 855.6923
$ws.Range("L97").Value = 2567.0769
$ws.Range("N97").Value = -3559.0769

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 2235.5
$ws.Range("I98").Value = 614
$ws.Range("K98").Value = 1842
$ws.Range("M98").Value = -344

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 31250836
$ws.Range("I107").Value = 166667000
$ws.Range("J107").Value = 951.53845
$ws.Range("K107").Value = 500001000
$ws.Range("L107").Value = 2854.61535
$ws.Range("M107").Value = -499999080
$ws.Range("N107").Value = -6694.61535

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3201
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 3445.5557
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 31010.0013
$ws.Range("M132").Value = -6470
$ws.Range("N132").Value = -36070.0013

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2167065
$ws.Range("I102").Value = 3573403.5
$ws.Range("J102").Value = 3467.5386
$ws.Range("K102").Value = 3573403.5
$ws.Range("L102").Value = 3467.5386
$ws.Range("M102").Value = -3571781.5
$ws.Range("N102").Value = -6711.5386

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4115.385
$ws.Range("I132").Value = 3653
$ws.Range("J132").Value = 4320.8887
$ws.Range("K132").Value = 10959
$ws.Range("L132").Value = 12962.6661
$ws.Range("M132").Value = -8429
$ws.Range("N132").Value = -18022.6661

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 40003520
$ws.Range("I40").Value = 76925370
$ws.Range("J40").Value = 4846.75
$ws.Range("K40").Value = 76925370
$ws.Range("L40").Value = 4846.75
$ws.Range("M40").Value = -76925234
$ws.Range("N40").Value = -5118.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6224.5835
$ws.Range("I61").Value = 5216.3335
$ws.Range("J61").Value = 7232.8335
$ws.Range("K61").Value = 5216.3335
$ws.Range("L61").Value = 7232.8335
$ws.Range("M61").Value = -5014.3335
$ws.Range("N61").Value = -7636.8335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 6224.5835
$ws.Range("I113").Value = 5216.3335
$ws.Range("J113").Value = 7232.8335
$ws.Range("K113").Value = 5216.3335
$ws.Range("L113").Value = 7232.8335
$ws.Range("M113").Value = -3046.3335
$ws.Range("N113").Value = -11572.8335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 32311.666
$ws.Range("J127").Value = 32311.666
$ws.Range("L127").Value = 32311.666
$ws.Range("N127").Value = -42231.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 3367.3333
$ws.Range("J101").Value = 3367.3333
$ws.Range("L101").Value = 3367.3333
$ws.Range("N101").Value = -9857.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H130").Value = 29766.924
$ws.Range("J130").Value = 29766.924
$ws.Range("L130").Value = 29766.924
$ws.Range("N130").Value = -39806.924
